# New added script and Updated Commit = 29/09/2020
#
# Walks through the "CriteriaLogic" city-value sheets, mirroring the user's
# click-through in Excel: select C2 on a couple of unrelated sheets, then
# update the city dropdown value on two sheets (new cities "Pratapgarh" and
# "Hydrabad"), and finish on the "Lead_Config29_2_CL" tab with C2 selected
# (making it the active tab instead of "Lead_Config21_CL").

$wb = $excel.ActiveWorkbook

# Lead_Config26_CL - just visited / C2 selected, value left as-is ("Noida").
$ws = $wb.Worksheets.Item("Lead_Config26_CL")
$ws.Activate()
$ws.Range("C2").Select()

# Lead_Config27_1_CL - city value changed to the newly added "Pratapgarh".
$ws = $wb.Worksheets.Item("Lead_Config27_1_CL")
$ws.Activate()
$ws.Range("C2").Value = "Pratapgarh"
$ws.Range("C2").Select()

# Lead_Config28_CL - just visited / C2 selected, value left as-is ("Hot").
$ws = $wb.Worksheets.Item("Lead_Config28_CL")
$ws.Activate()
$ws.Range("C2").Select()

# Lead_Config29_1_CL - city value changed to the newly added "Hydrabad".
$ws = $wb.Worksheets.Item("Lead_Config29_1_CL")
$ws.Activate()
$ws.Range("C2").Value = "Hydrabad"
$ws.Range("C2").Select()

# Lead_Config29_2_CL - ends up as the active/selected tab with C2 selected
# (previously the selection was on A2, and "Lead_Config21_CL" was the tab
# left active).
$ws = $wb.Worksheets.Item("Lead_Config29_2_CL")
$ws.Activate()
$ws.Range("C2").Select()
